$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "302.92"
Set-TextValue "E2" "2.41%"
Set-TextValue "G2" "7"
Set-TextValue "D3" "31.71"
Set-TextValue "E3" "0.50%"
Set-TextValue "G3" "7"
Set-TextValue "D4" "5.165"
Set-TextValue "E4" "1.13%"
Set-TextValue "G4" "7"
Set-TextValue "D5" "0.07817"
Set-TextValue "E5" "4.13%"
Set-TextValue "G5" "7"
Set-TextValue "D6" "2.441"
Set-TextValue "E6" "43.70%"
Set-TextValue "G6" "7"
Set-TextValue "D7" "7.970"
Set-TextValue "E7" "3.11%"
Set-TextValue "G7" "7"
Set-TextValue "D8" "3.872"
Set-TextValue "E8" "1.87%"
Set-TextValue "G8" "7"
Set-TextValue "D9" "0.9114"
Set-TextValue "E9" "-2.37%"
Set-TextValue "G9" "7"
Set-TextValue "D10" "0.1729"
Set-TextValue "E10" "2.35%"
Set-TextValue "G10" "7"
Set-TextValue "D11" "0.07336"
Set-TextValue "E11" "-1.20%"
Set-TextValue "G11" "7"
Set-TextValue "D12" "0.08147"
Set-TextValue "E12" "2.40%"
Set-TextValue "G12" "7"
Set-TextValue "D13" "0.03043"
Set-TextValue "E13" "0.65%"
Set-TextValue "G13" "7"
Set-TextValue "D14" "0.09941"
Set-TextValue "E14" "0.44%"
Set-TextValue "G14" "7"
Set-TextValue "D15" "0.001518"
Set-TextValue "E15" "1.16%"
Set-TextValue "G15" "7"
Set-TextValue "D16" "0.006010"
Set-TextValue "E16" "-5.47%"
Set-TextValue "G16" "7"
Set-TextValue "D17" "3.497"
Set-TextValue "E17" "1.08%"
Set-TextValue "G17" "7"
Set-TextValue "E18" "1.09%"
Set-TextValue "G18" "7"
Set-TextValue "D19" "0.3243"
Set-TextValue "E19" "-1.10%"
Set-TextValue "G19" "7"
Set-TextValue "D20" "0.1338"
Set-TextValue "E20" "0.78%"
Set-TextValue "G20" "7"
Set-TextValue "D21" "4.692"
Set-TextValue "E21" "2.78%"
Set-TextValue "G21" "7"
Set-TextValue "D22" "0.04655"
Set-TextValue "E22" "0.14%"
Set-TextValue "G22" "7"
Set-TextValue "D23" "0.1565"
Set-TextValue "E23" "0.30%"
Set-TextValue "G23" "7"
Set-TextValue "D24" "0.001261"
Set-TextValue "E24" "3.49%"
Set-TextValue "G24" "7"
Set-TextValue "D25" "0.004518"
Set-TextValue "E25" "2.27%"
Set-TextValue "G25" "7"
Set-TextValue "E26" "3.64%"
Set-TextValue "G26" "7"
Set-TextValue "D27" "0.0002742"
Set-TextValue "E27" "45.95%"
Set-TextValue "G27" "7"
Set-TextValue "G28" "7"
Set-TextValue "G29" "7"
Set-TextValue "G30" "7"
Set-TextValue "G31" "7"
Set-TextValue "G32" "7"
Set-TextValue "G33" "7"
Set-TextValue "G34" "7"
Set-TextValue "G35" "7"
Set-TextValue "G36" "7"
Set-TextValue "G37" "7"
Set-TextValue "G38" "7"
Set-TextValue "D39" "0.01779"
Set-TextValue "E39" "7.27%"
Set-TextValue "G39" "7"
Set-TextValue "D40" "0.04563"
Set-TextValue "E40" "2.21%"
Set-TextValue "G40" "7"
Set-TextValue "D41" "0.007306"
Set-TextValue "E41" "3.40%"
Set-TextValue "G41" "7"
Set-TextValue "D42" "0.1361"
Set-TextValue "E42" "2.62%"
Set-TextValue "G42" "7"
Set-TextValue "D43" "0.002240"
Set-TextValue "E43" "8.56%"
Set-TextValue "G43" "7"
Set-TextValue "D44" "0.01076"
Set-TextValue "E44" "-4.15%"
Set-TextValue "G44" "7"
Set-TextValue "D45" "0.00006482"
Set-TextValue "E45" "8.07%"
Set-TextValue "G45" "7"
Set-TextValue "E46" "-0.03%"
Set-TextValue "G46" "7"
Set-TextValue "G47" "7"
Set-TextValue "D48" "0.009898"
Set-TextValue "E48" "-23.66%"
Set-TextValue "G48" "7"
Set-TextValue "D49" "0.00002100"
Set-TextValue "E49" "-0.03%"
Set-TextValue "G49" "7"
Set-TextValue "D50" "0.0002000"
Set-TextValue "E50" "0.04%"
Set-TextValue "G50" "7"
Set-TextValue "G51" "7"
